$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Experimental flag: true -> false (must remain a text cell, not boolean)
$ws.Range("B7").Value = "'false"

# Date updated
$ws.Range("B8").Value = "2025-10-03T16:37:46+01:00"
